$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Relocate the "_GoBack" bookmark.
#    Originally it wraps the "Anforderungsanalyse" bullet near the end of the
#    document (a leftover from the author's last edit there). After exporting
#    to PDF, Word re-records "_GoBack" at the position of the most recent
#    cursor/edit, which in this revision is the very start of the document
#    (the title). So: delete the old one, create a fresh zero-length bookmark
#    at the top of the document.
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# A collapsed bookmark placed exactly at absolute position 0 tends to swallow
# the whole first paragraph when serialized, so: temporarily insert a single
# placeholder character at position 0, drop the (now non-boundary-0) bookmark
# right after it, then remove the placeholder again. The bookmark stays put
# as a proper zero-length bookmark at the true start of the document.
$placeholder = $d.Range(0, 0)
$placeholder.InsertBefore("X")

$bmRange = $d.Range(1, 1)
$d.Bookmarks.Add("_GoBack", $bmRange)

$placeholder2 = $d.Range(0, 1)
$placeholder2.Text = ""

# ---------------------------------------------------------------------------
# 2) The footer's "Seite <n>" PAGE field has a stale cached result ("3")
#    from before the edit; refresh its displayed text to "2" to match the
#    recalculated page count after export.
# ---------------------------------------------------------------------------
$footer = $d.Sections.Item(1).Footers.Item(1)
$footerRange = $footer.Range
$footerText = $footerRange.Text
$label = "Seite "
$labelIdx = $footerText.IndexOf($label)
if ($labelIdx -ge 0) {
    $charPos = $labelIdx + $label.Length + 1
    $pageNoChar = $footerRange.Characters.Item($charPos)
    if ($pageNoChar.Text -eq "3") {
        $pageNoChar.Text = "2"
    }
}
